# Update "想去人数" (want-to-go count) values on both the "展览" sheet and
# the "全部类型" sheet (which mirrors the same data).
#
#   Sheet "展览"   F3: 180 -> 182
#   Sheet "展览"   F4: 133 -> 134
#   Sheet "全部类型" F3: 180 -> 182
#   Sheet "全部类型" F4: 133 -> 134

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 182
    $ws.Range("F4").Value = 134
}
